$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record is inserted before the existing row 47, pushing the
# old rows 47 and 48 down to rows 48 and 49 respectively (their contents
# stay the same). The new row 47 holds the new "Black Amber" observation.
$ws.Rows("47:47").Insert()

$ws.Cells.Item(47, 1).Value = 11
$ws.Cells.Item(47, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(47, 3).Value = "Bíobío"
$ws.Cells.Item(47, 4).Value = 44595
$ws.Cells.Item(47, 5).Value = 8
$ws.Cells.Item(47, 6).Value = "Fruta"
$ws.Cells.Item(47, 7).Value = 100103
$ws.Cells.Item(47, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(47, 9).Value = 100103002
$ws.Cells.Item(47, 10).Value = "Ciruela"
$ws.Cells.Item(47, 11).Value = "Black Amber"
$ws.Cells.Item(47, 12).Value = "Primera"
$ws.Cells.Item(47, 13).Value = 250
$ws.Cells.Item(47, 14).Value = 8500
$ws.Cells.Item(47, 15).Value = 9000
$ws.Cells.Item(47, 16).Value = 8740
$ws.Cells.Item(47, 17).Value = "`$/caja 16 kilos granel"
$ws.Cells.Item(47, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(47, 19).Value = 546
$ws.Cells.Item(47, 20).Value = 16

# Keep the date column's existing number format (same style the rest of
# the "Fecha" column already uses) on the freshly inserted cell.
$ws.Cells.Item(47, 4).NumberFormat = $ws.Cells.Item(48, 4).NumberFormat
